# Add migraine information to create conceptsets
# Adds three new rows (14-16) to Sheet 1, mirroring the existing
# GDM_diagnoses / PE_diagnoses blocks but for a new "Migraine_diagnoses"
# concept set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - EVENTS
$ws.Range("A14").Value = "EVENTS"
$ws.Range("B14").Value = "ARS"
$ws.Range("C14").Value = "Migraine_diagnoses"
$ws.Range("D14").Value = "codesheet"
$ws.Range("E14").Value = "code"
$ws.Range("F14").Value = "vocabulary"
$ws.Range("J14").Value = "event_code"
$ws.Range("K14").Value = "event_record_vocabulary"
$ws.Range("O14").Value = "start_date_record"

# Row 15 - MEDICAL_OBSERVATIONS
$ws.Range("A15").Value = "MEDICAL_OBSERVATIONS"
$ws.Range("B15").Value = "ARS"
$ws.Range("C15").Value = "Migraine_diagnoses"
$ws.Range("D15").Value = "codesheet"
$ws.Range("E15").Value = "code"
$ws.Range("F15").Value = "vocabulary"
$ws.Range("J15").Value = "mo_code"
$ws.Range("K15").Value = "mo_record_vocabulary"
$ws.Range("O15").Value = "mo_date"

# Row 16 - SURVEY_OBSERVATIONS
$ws.Range("A16").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B16").Value = "ARS"
$ws.Range("C16").Value = "Migraine_diagnoses"
$ws.Range("D16").Value = "codesheet"
$ws.Range("E16").Value = "code"
$ws.Range("F16").Value = "vocabulary"
$ws.Range("J16").Value = "so_source_value"
$ws.Range("K16").Value = "so_unit"
$ws.Range("O16").Value = "so_date"

# Update the selection to match the author's final cursor position
$ws.Range("C20").Select()
